$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new values to row 4 (B4:D4) - continues the SLA row started at A4
# (order matters for shared-string table insertion order)
$ws.Range("C4").Value = "CAAssignService"
$ws.Range("D4").Value = "PegaCS-Cases:08-06-01"
$ws.Range("B4").Value = "PegaCA-Work-Service-GeneralRequest"

# Apply left/top alignment style to D4 (new cellXfs entry)
$ws.Range("D4").HorizontalAlignment = -4131
$ws.Range("D4").VerticalAlignment = -4160

# Widen column D to fit the new (longer) content, mirroring Excel's
# "best fit" auto-resize that happens after typing into the column
$ws.Columns.Item(4).ColumnWidth = 20.33

# Update the active selection / view
$ws.Range("B4").Select()
